$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 0.04000434785130035
$ws.Range("E2").Value = 0.03950676534389326
$ws.Range("F2").Value = 0.03939952808019245
$ws.Range("G2").Value = 0.03933405950324943
$ws.Range("D3").Value = 0.04079397860326064
$ws.Range("E3").Value = 0.04068034854710206
$ws.Range("F3").Value = 0.04020477829632312
$ws.Range("G3").Value = 0.03943734061124522
$ws.Range("D4").Value = 0.03277589664240757
$ws.Range("E4").Value = 0.03266872209705878
$ws.Range("F4").Value = 0.03248301738290605
$ws.Range("G4").Value = 0.03244428900374872
$ws.Range("C5").Value = 0.03656281317052121
$ws.Range("D5").Value = 0.03623691955667992
$ws.Range("E5").Value = 0.03585351322284168
$ws.Range("F5").Value = 0.0355985663950515
$ws.Range("G5").Value = 0.0355226383124756
$ws.Range("D6").Value = 0.03728456689129077
$ws.Range("F6").Value = 0.03707158257353536
$ws.Range("G6").Value = 0.03705329114574955
$ws.Range("F7").Value = 0.0406940033749759
$ws.Range("G7").Value = 0.0405234318667373
$ws.Range("C8").Value = 0.03421651979822254
$ws.Range("D8").Value = 0.03419497371930777
$ws.Range("G8").Value = 0.03366394205114161
$ws.Range("E9").Value = 0.03560290821817723
$ws.Range("F9").Value = 0.03531960457970115
$ws.Range("G9").Value = 0.03495653006861789
$ws.Range("C10").Value = 0.03653259067730615
$ws.Range("D10").Value = 0.03552617096385219
$ws.Range("F10").Value = 0.03519423961164765
$ws.Range("G10").Value = 0.0350235436339904
$ws.Range("C11").Value = 0.04091767521764635
$ws.Range("D11").Value = 0.0385188743623843
$ws.Range("E11").Value = 0.03844635294504702
$ws.Range("F11").Value = 0.03834185834755878
$ws.Range("G11").Value = 0.03819325945150712
$ws.Range("D12").Value = 0.03792180142592391
$ws.Range("E12").Value = 0.03742173350191724
$ws.Range("F12").Value = 0.03715957935338052
$ws.Range("G12").Value = 0.03706082751841745
$ws.Range("C13").Value = 0.03338453141001082
$ws.Range("D13").Value = 0.03313271733065908
$ws.Range("E13").Value = 0.03270756990864491
$ws.Range("F13").Value = 0.03258181135511897
$ws.Range("G13").Value = 0.03257522049311584
$ws.Range("D14").Value = 0.03360184159643557
$ws.Range("E14").Value = 0.03343128430985657
$ws.Range("F14").Value = 0.03335374440899627
$ws.Range("G14").Value = 0.03328665432776195
$ws.Range("C15").Value = 0.03539057302831939
$ws.Range("D15").Value = 0.03537070315239201
$ws.Range("E15").Value = 0.03498744106857717
$ws.Range("G15").Value = 0.0346362429200126
$ws.Range("C16").Value = 0.03991546163275893
$ws.Range("D16").Value = 0.03918820274350299
$ws.Range("E16").Value = 0.03675132126988168
$ws.Range("F16").Value = 0.0366958213108294
$ws.Range("G16").Value = 0.0360152789268874
$ws.Range("E17").Value = 0.03566394302020694
$ws.Range("G17").Value = 0.03448518566695776
$ws.Range("D18").Value = 0.03566312548035062
$ws.Range("E18").Value = 0.03471337128226992
$ws.Range("F18").Value = 0.03433844397141863
$ws.Range("G18").Value = 0.03428761527343555
$ws.Range("C19").Value = 0.03095735838686689
$ws.Range("D19").Value = 0.0303985255414433
$ws.Range("E19").Value = 0.0303829120035985
$ws.Range("C20").Value = 0.03466956171178111
$ws.Range("D20").Value = 0.03429314801109391
$ws.Range("E20").Value = 0.03428316246761137
$ws.Range("C21").Value = 0.03539235174792499
$ws.Range("E21").Value = 0.03456261560107315
$ws.Range("F21").Value = 0.03435550735347877
$ws.Range("G21").Value = 0.03429948596094831
$ws.Range("C22").Value = 0.03462154954958127
$ws.Range("D22").Value = 0.03336659443586466
$ws.Range("E22").Value = 0.03259913877538174
$ws.Range("C23").Value = 0.03598883961844532
$ws.Range("D23").Value = 0.03550380072172239
$ws.Range("E23").Value = 0.03537222967315393
$ws.Range("F23").Value = 0.03464846926148695
$ws.Range("G23").Value = 0.03446625873392323
$ws.Range("C24").Value = 0.03280995927863482
$ws.Range("D24").Value = 0.03256896951684068
$ws.Range("E24").Value = 0.03236041572849568
$ws.Range("D25").Value = 0.033299072809871
$ws.Range("E25").Value = 0.03297575741305355
$ws.Range("G25").Value = 0.03282084132458501
$ws.Range("D26").Value = 0.03695378827398092
$ws.Range("E26").Value = 0.03479153338419285
$ws.Range("F26").Value = 0.03462708721544569
$ws.Range("C27").Value = 0.0323376042408559
$ws.Range("D27").Value = 0.03002414028355591
$ws.Range("E27").Value = 0.03000353857958505
$ws.Range("F27").Value = 0.02990512759002684
$ws.Range("G27").Value = 0.0298542206808492
$ws.Range("D28").Value = 0.03324318109941769
$ws.Range("E28").Value = 0.03293092705757739
$ws.Range("F28").Value = 0.03257011401562933
$ws.Range("C29").Value = 0.03138532251684729
$ws.Range("D29").Value = 0.03098007193113086
$ws.Range("E29").Value = 0.03071413927582397
$ws.Range("C30").Value = 0.03638422517434112
$ws.Range("D30").Value = 0.03576772347051604
$ws.Range("E30").Value = 0.03528182982354162
$ws.Range("F30").Value = 0.03516716498859982
$ws.Range("G30").Value = 0.03515581262048962
$ws.Range("D31").Value = 0.03647959022505102
$ws.Range("E31").Value = 0.03639609638579747
$ws.Range("F31").Value = 0.03608925287962519
$ws.Range("G31").Value = 0.03585238040580251
$ws.Range("C32").Value = 0.0297001642221492
$ws.Range("D32").Value = 0.02941747981201072
$ws.Range("E32").Value = 0.02928659114696261
$ws.Range("F32").Value = 0.02926592739347359
$ws.Range("C33").Value = 0.03327279260527567
$ws.Range("D33").Value = 0.03296669564690752
$ws.Range("E33").Value = 0.03263844303355716
$ws.Range("F33").Value = 0.03234932350551908
$ws.Range("C34").Value = 0.04212884937159443
$ws.Range("D34").Value = 0.04118275274540119
$ws.Range("E34").Value = 0.03762351655090383
$ws.Range("F34").Value = 0.03741446571238259
$ws.Range("G34").Value = 0.03714316980486784
$ws.Range("C35").Value = 0.03022136572669462
$ws.Range("D35").Value = 0.02988699167425481
$ws.Range("E35").Value = 0.02986337046042326
$ws.Range("G35").Value = 0.02944383654311555
$ws.Range("F36").Value = 0.03393920227359754
$ws.Range("G36").Value = 0.03310555274254069
$ws.Range("C37").Value = 0.03457403334361462
$ws.Range("D37").Value = 0.03406767463599975
$ws.Range("E37").Value = 0.03349918029503234
$ws.Range("F37").Value = 0.03331448502224645
$ws.Range("G37").Value = 0.03309613849809422
$ws.Range("C38").Value = 0.02413056364637791
$ws.Range("D38").Value = 0.02351458093364857
$ws.Range("E38").Value = 0.02340915108963856
$ws.Range("F38").Value = 0.02329687358705899
$ws.Range("G38").Value = 0.02328501571895603
$ws.Range("C39").Value = 0.03141365961958499
$ws.Range("D39").Value = 0.03130418992312569
$ws.Range("C40").Value = 0.03263516948121871
$ws.Range("D40").Value = 0.0321085183084026
$ws.Range("E40").Value = 0.03167352412025241
$ws.Range("F40").Value = 0.03158643039525289
$ws.Range("G40").Value = 0.03142402970055166
$ws.Range("C41").Value = 0.02635329007205709
$ws.Range("D41").Value = 0.02519793849179802
$ws.Range("E41").Value = 0.02514114134698422
$ws.Range("G41").Value = 0.02481891286673979
$ws.Range("C42").Value = 0.03477247789850443
$ws.Range("D42").Value = 0.03433907997604539
$ws.Range("F42").Value = 0.03261695793537692
$ws.Range("G42").Value = 0.03056010536682863
$ws.Range("C43").Value = 0.03598029452350893
$ws.Range("D43").Value = 0.03535989878123726
$ws.Range("E43").Value = 0.03489846940228675
$ws.Range("F43").Value = 0.03485106074628356
$ws.Range("G43").Value = 0.0347359466192077
$ws.Range("D44").Value = 0.03619748282849886
$ws.Range("E44").Value = 0.03607548577692608
$ws.Range("F44").Value = 0.03600204971800382
$ws.Range("G44").Value = 0.03560675485257794
$ws.Range("C45").Value = 0.0333266753660377
$ws.Range("D45").Value = 0.03280877799940893
$ws.Range("E45").Value = 0.0327195548569437
$ws.Range("F45").Value = 0.03199670698720346
$ws.Range("G45").Value = 0.03113721245805416
$ws.Range("D46").Value = 0.03269917374121813
$ws.Range("E46").Value = 0.03206698109689552
$ws.Range("F46").Value = 0.03199116658642889
$ws.Range("G46").Value = 0.03194744999059786
$ws.Range("D47").Value = 0.02124313864217251
$ws.Range("F47").Value = 0.02104881108638292
$ws.Range("C48").Value = 0.02896520685920233
$ws.Range("D48").Value = 0.02894123505625453
$ws.Range("E48").Value = 0.02791932962396472
$ws.Range("G48").Value = 0.02645542788022004
$ws.Range("D49").Value = 0.03045832944362489
$ws.Range("G49").Value = 0.02825516478047194
$ws.Range("D50").Value = 0.02478392804147314
$ws.Range("E50").Value = 0.02464927472267159
$ws.Range("F50").Value = 0.02446196962891726
$ws.Range("G50").Value = 0.02434551346842735
